# Update dashboards - 2025-10-28
# Refresh the FRED-sourced "Latest Date" and Present/Lag1-4 readings for the
# rows whose underlying series rolled forward this week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to round-trip a cell's original number format (so that
# writing a date-looking string into column N does not get auto-converted
# into an Excel date serial / auto-assigned a new "Text" style). We stash the
# target cell's current format onto the scratch cell, overwrite the value
# (with the cell temporarily forced to Text so the string isn't re-typed as a
# date), then paste the stashed format back on top so the style index is
# exactly what it was before.
$scratch = $ws.Range("AA200")

function Set-LatestDate($cellRef, $newDate) {
    $target = $ws.Range($cellRef)
    $target.Copy()
    $scratch.PasteSpecial(-4122)   # xlPasteFormats -> stash current format
    $target.NumberFormat = "@"     # force Text so the date string isn't re-typed
    $target.Value = $newDate
    $scratch.Copy()
    $target.PasteSpecial(-4122)    # xlPasteFormats -> restore original format/style
    $scratch.Clear()
}

# --- Row 29: T5YIFR (5yr, 5yr Forward) ---
Set-LatestDate "N29" "2025-10-27"
$ws.Range("Q29").Value = 2.19

# --- Row 30: T10YIE (10yr TIPS) ---
Set-LatestDate "N30" "2025-10-27"
$ws.Range("Q30").Value = 2.28
$ws.Range("R30").Value = 2.29
$ws.Range("S30").Value = 2.3
$ws.Range("T30").Value = 2.29
$ws.Range("U30").Value = 2.28

# --- Row 39: DTWEXBGS level (Nominal Broad US Dollar Index) ---
Set-LatestDate "N39" "2025-10-24"
$ws.Range("Q39").Value = 121.342
$ws.Range("R39").Value = 121.3633
$ws.Range("S39").Value = 121.3075
$ws.Range("T39").Value = 121.302
$ws.Range("U39").Value = 121.0394

# --- Row 40: DTWEXBGS Y/Y % Delta ---
Set-LatestDate "N40" "2025-10-24"
$ws.Range("Q40").Value = -0.02945577953796328

# --- Row 48: DGS2 (2y UST) ---
Set-LatestDate "N48" "2025-10-24"
$ws.Range("R48").Value = 3.48
$ws.Range("T48").Value = 3.45

# --- Row 49: DGS5 (5y UST) ---
Set-LatestDate "N49" "2025-10-24"
$ws.Range("R49").Value = 3.61
$ws.Range("T49").Value = 3.56
$ws.Range("U49").Value = 3.58

# --- Row 50: DGS10 (10y UST) ---
Set-LatestDate "N50" "2025-10-24"
$ws.Range("Q50").Value = 4.02
$ws.Range("R50").Value = 4.01
$ws.Range("S50").Value = 3.97
$ws.Range("T50").Value = 3.98
$ws.Range("U50").Value = 4

# --- Row 52: DBAA (BAA) ---
Set-LatestDate "N52" "2025-10-24"
$ws.Range("R52").Value = 5.67
$ws.Range("S52").Value = 5.66
$ws.Range("T52").Value = 5.65
$ws.Range("U52").Value = 5.68
